$p = $ppt.ActivePresentation

# --- Slide 3 ("references"): Content Placeholder 3 -----------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)

# Widen the placeholder (EMU 6832600 -> 7787640, i.e. pts 538.0 -> 613.2)
$sh3.Width = 613.2

$tr3 = $sh3.TextFrame.TextRange

# Append a trailing space to the existing sentence (stays in the same run)
$run3a = $tr3.Find("Is an alias for something else.", 0)
$run3a.Text = "Is an alias for something else. "

# Append the new red code snippet as its own run right after it
$para3 = $tr3.Paragraphs(1, 1)
$origLen3 = $para3.Text.TrimEnd("`r").Length
$para3.InsertAfter("Vehicle& v2 = v1;") | Out-Null
$run3b = $tr3.Characters($origLen3 + 1, 17)
$run3b.Font.Color.RGB = 255

# --- Slide 4 ("Pointers"): Content Placeholder 3 --------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)

# Widen the placeholder (EMU 6832600 -> 8879840, i.e. pts 538.0 -> 699.2)
$sh4.Width = 699.2

$tr4 = $sh4.TextFrame.TextRange

# Append a trailing space to the existing sentence (stays in the same run)
$run4a = $tr4.Find("An alias to something else.", 0)
$run4a.Text = "An alias to something else. "

$para4 = $tr4.Paragraphs(1, 1)
$pos4 = $para4.Text.TrimEnd("`r").Length

# "Vehicle* " in red
$para4.InsertAfter("Vehicle* ") | Out-Null
$run4b = $tr4.Characters($pos4 + 1, 9)
$run4b.Font.Color.RGB = 255
$pos4 = $pos4 + 9

# "vp" in red
$para4.InsertAfter("vp") | Out-Null
$run4c = $tr4.Characters($pos4 + 1, 2)
$run4c.Font.Color.RGB = 255
$pos4 = $pos4 + 2

# " = &v1;" in red
$para4.InsertAfter(" = &v1;") | Out-Null
$run4d = $tr4.Characters($pos4 + 1, 7)
$run4d.Font.Color.RGB = 255
